$d = $word.ActiveDocument

function ReplaceText($find, $replace) {
    $range = $d.Content
    $ok = $range.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Host "NOT FOUND: $find"
    }
}

# Small line-number tweaks in the stack trace
ReplaceText "JavaMethodService.java:163" "JavaMethodService.java:162"
ReplaceText "AbstractService.java:136" "AbstractService.java:135"
ReplaceText "EvaluationServices.java:168" "EvaluationServices.java:172"
ReplaceText "AstEvaluator.java:189" "AstEvaluator.java:186"
ReplaceText "AstSwitch.java:118" "AstSwitch.java:119"
ReplaceText "AstEvaluator.java:112" "AstEvaluator.java:109"
ReplaceText "GeneratedMethodAccessor74" "GeneratedMethodAccessor73"

# Replace the long tail of the stack trace (surefire/equinox/eclipse launcher frames)
# with the shorter jdt-based frames.
$oldTail = "	at org.apache.maven.surefire.junit4.JUnit4Provider.execute(JUnit4Provider.java:264)" + [char]10 + `
"	at org.apache.maven.surefire.junit4.JUnit4Provider.executeTestSet(JUnit4Provider.java:153)" + [char]10 + `
"	at org.apache.maven.surefire.junit4.JUnit4Provider.invoke(JUnit4Provider.java:124)" + [char]10 + `
"	at sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)" + [char]10 + `
"	at sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)" + [char]10 + `
"	at sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)" + [char]10 + `
"	at java.lang.reflect.Method.invoke(Method.java:498)" + [char]10 + `
"	at org.apache.maven.surefire.util.ReflectionUtils.invokeMethodWithArray2(ReflectionUtils.java:208)" + [char]10 + `
"	at org.apache.maven.surefire.booter.ProviderFactory`$ProviderProxy.invoke(ProviderFactory.java:156)" + [char]10 + `
"	at org.apache.maven.surefire.booter.ProviderFactory.invokeProvider(ProviderFactory.java:82)" + [char]10 + `
"	at org.eclipse.tycho.surefire.osgibooter.OsgiSurefireBooter.run(OsgiSurefireBooter.java:91)" + [char]10 + `
"	at org.eclipse.tycho.surefire.osgibooter.HeadlessTestApplication.run(HeadlessTestApplication.java:21)" + [char]10 + `
"	at sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)" + [char]10 + `
"	at sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)" + [char]10 + `
"	at sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)" + [char]10 + `
"	at java.lang.reflect.Method.invoke(Method.java:498)" + [char]10 + `
"	at org.eclipse.equinox.internal.app.EclipseAppContainer.callMethodWithException(EclipseAppContainer.java:587)" + [char]10 + `
"	at org.eclipse.equinox.internal.app.EclipseAppHandle.run(EclipseAppHandle.java:198)" + [char]10 + `
"	at org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.runApplication(EclipseAppLauncher.java:134)" + [char]10 + `
"	at org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.start(EclipseAppLauncher.java:104)" + [char]10 + `
"	at org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:388)" + [char]10 + `
"	at org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:243)" + [char]10 + `
"	at sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)" + [char]10 + `
"	at sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)" + [char]10 + `
"	at sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)" + [char]10 + `
"	at java.lang.reflect.Method.invoke(Method.java:498)" + [char]10 + `
"	at org.eclipse.equinox.launcher.Main.invokeFramework(Main.java:656)" + [char]10 + `
"	at org.eclipse.equinox.launcher.Main.basicRun(Main.java:592)" + [char]10 + `
"	at org.eclipse.equinox.launcher.Main.run(Main.java:1498)" + [char]10 + `
"	at org.eclipse.equinox.launcher.Main.main(Main.java:1471)"

$newTail = "	at org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)" + [char]10 + `
"	at org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)" + [char]10 + `
"	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538)" + [char]10 + `
"	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:760)" + [char]10 + `
"	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:460)" + [char]10 + `
"	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:206)"

$range2 = $d.Content
$ok2 = $range2.Find.Execute($oldTail, $false, $true, $false, $false, $false, $true, 1, $false, $newTail, 2)
if (-not $ok2) {
    Write-Host "NOT FOUND: tail block"
} else {
    Write-Host "Tail block replaced"
}
